$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7100665
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45172.58333333334
$arr[0,3] = 'OFI Crete'
$arr[0,4] = 'PAOK Salonika'
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 4.75
$arr[0,9] = 3.6
$arr[0,10] = 1.75
$arr[0,11] = 4.5
$arr[0,12] = 3.8
$arr[0,13] = 1.75
$arr[0,14] = 0.75
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2.75
$arr[0,18] = 2
$arr[0,19] = 1.85
$arr[0,20] = 3.5
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8500000000000001
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.8500000000000001
$ws.Range("B15:AB15").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7100664
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45172.58333333334
$arr[0,3] = 'Olympiakos'
$arr[0,4] = 'Lamia'
$arr[0,5] = 4
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.25
$arr[0,9] = 5.5
$arr[0,10] = 12
$arr[0,11] = 1.222
$arr[0,12] = 6
$arr[0,13] = 15
$arr[0,14] = -1.75
$arr[0,15] = 1.875
$arr[0,16] = 1.975
$arr[0,17] = 2.75
$arr[0,18] = 1.825
$arr[0,19] = 2.025
$arr[0,20] = 0.222
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.875
$arr[0,24] = -1
$arr[0,25] = 0.825
$arr[0,26] = -1
$ws.Range("B16:AB16").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7100661
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45172.58333333334
$arr[0,3] = 'Aris Salonika'
$arr[0,4] = 'Asteras Tripolis'
$arr[0,5] = 3
$arr[0,6] = 2
$arr[0,7] = 'H'
$arr[0,8] = 1.8
$arr[0,9] = 3.4
$arr[0,10] = 4.75
$arr[0,11] = 1.55
$arr[0,12] = 3.8
$arr[0,13] = 7
$arr[0,14] = -1
$arr[0,15] = 2
$arr[0,16] = 1.85
$arr[0,17] = 2.25
$arr[0,18] = 1.825
$arr[0,19] = 2.025
$arr[0,20] = 0.55
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0.825
$arr[0,26] = -1
$ws.Range("B17:AB17").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937185
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45196.60416666666
$arr[0,3] = 'Olympiakos'
$arr[0,4] = 'Aris Salonika'
$arr[0,5] = 4
$arr[0,6] = 1
$arr[0,7] = 'H'
$arr[0,8] = 1.4
$arr[0,9] = 4.333
$arr[0,10] = 8.5
$arr[0,11] = 1.4
$arr[0,12] = 4.2
$arr[0,13] = 9
$arr[0,14] = -1.25
$arr[0,15] = 2
$arr[0,16] = 1.85
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.875
$arr[0,20] = 0.3999999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1
$arr[0,24] = -1
$arr[0,25] = 0.9750000000000001
$arr[0,26] = -1
$ws.Range("B36:AB36").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937187
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45196.60416666666
$arr[0,3] = 'PAOK Salonika'
$arr[0,4] = 'Volos NFC'
$arr[0,5] = 3
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.166
$arr[0,9] = 7.5
$arr[0,10] = 13
$arr[0,11] = 1.2
$arr[0,12] = 7.5
$arr[0,13] = 11
$arr[0,14] = -2
$arr[0,15] = 1.95
$arr[0,16] = 1.9
$arr[0,17] = 3.25
$arr[0,18] = 1.975
$arr[0,19] = 1.875
$arr[0,20] = 0.2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.95
$arr[0,24] = -1
$arr[0,25] = -0.5
$arr[0,26] = 0.4375
$ws.Range("B37:AB37").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937192
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45221.41666666666
$arr[0,3] = 'Asteras Tripolis'
$arr[0,4] = 'AEK Athens'
$arr[0,5] = 0
$arr[0,6] = 3
$arr[0,7] = 'A'
$arr[0,8] = 6.5
$arr[0,9] = 3.6
$arr[0,10] = 1.6
$arr[0,11] = 7.5
$arr[0,12] = 4
$arr[0,13] = 1.5
$arr[0,14] = 1
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.5
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.5
$arr[0,23] = -1
$arr[0,24] = 0.825
$arr[0,25] = 0.925
$arr[0,26] = -1
$ws.Range("B51:AB51").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6935736
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45221.41666666666
$arr[0,3] = 'Kifisias FC'
$arr[0,4] = 'OFI Crete'
$arr[0,5] = 0
$arr[0,6] = 0
$arr[0,7] = 'D'
$arr[0,8] = 3.3
$arr[0,9] = 3.25
$arr[0,10] = 2.25
$arr[0,11] = 3.1
$arr[0,12] = 3.4
$arr[0,13] = 2.4
$arr[0,14] = 0.25
$arr[0,15] = 1.8
$arr[0,16] = 2.05
$arr[0,17] = 2.5
$arr[0,18] = 1.825
$arr[0,19] = 2.025
$arr[0,20] = -1
$arr[0,21] = 2.4
$arr[0,22] = -1
$arr[0,23] = 0.4
$arr[0,24] = -0.5
$arr[0,25] = -1
$arr[0,26] = 1.025
$ws.Range("B52:AB52").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937238
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45305.64583333334
$arr[0,3] = 'PAOK Salonika'
$arr[0,4] = 'Giannina'
$arr[0,5] = 4
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.111
$arr[0,9] = 9
$arr[0,10] = 23
$arr[0,11] = 1.25
$arr[0,12] = 6
$arr[0,13] = 9
$arr[0,14] = -1.75
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.75
$arr[0,18] = 1.8
$arr[0,19] = 2.05
$arr[0,20] = 0.25
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1.025
$arr[0,24] = -1
$arr[0,25] = 0.8
$arr[0,26] = -1
$ws.Range("B124:AB124").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6936857
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45305.64583333334
$arr[0,3] = 'AEK Athens'
$arr[0,4] = 'Panathinaikos'
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = 'D'
$arr[0,8] = 1.909
$arr[0,9] = 3.5
$arr[0,10] = 4.2
$arr[0,11] = 2.15
$arr[0,12] = 3.2
$arr[0,13] = 3.5
$arr[0,14] = -0.25
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2
$arr[0,18] = 1.8
$arr[0,19] = 2.05
$arr[0,20] = -1
$arr[0,21] = 2.2
$arr[0,22] = -1
$arr[0,23] = -0.5
$arr[0,24] = 0.5
$arr[0,25] = 0.8
$arr[0,26] = -1
$ws.Range("B125:AB125").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937250
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45326.52083333334
$arr[0,3] = 'Giannina'
$arr[0,4] = 'Lamia'
$arr[0,5] = 1
$arr[0,6] = 4
$arr[0,7] = 'A'
$arr[0,8] = 2.3
$arr[0,9] = 3.25
$arr[0,10] = 3.25
$arr[0,11] = 2.55
$arr[0,12] = 2.875
$arr[0,13] = 3.1
$arr[0,14] = 0
$arr[0,15] = 1.75
$arr[0,16] = 2.125
$arr[0,17] = 2
$arr[0,18] = 1.85
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.1
$arr[0,23] = -1
$arr[0,24] = 1.125
$arr[0,25] = 0.8500000000000001
$arr[0,26] = -1
$ws.Range("B143:AB143").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937247
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45326.52083333334
$arr[0,3] = 'AEK Athens'
$arr[0,4] = 'Asteras Tripolis'
$arr[0,5] = 4
$arr[0,6] = 2
$arr[0,7] = 'H'
$arr[0,8] = 1.285
$arr[0,9] = 5.5
$arr[0,10] = 12
$arr[0,11] = 1.285
$arr[0,12] = 5.75
$arr[0,13] = 10
$arr[0,14] = -1.5
$arr[0,15] = 1.825
$arr[0,16] = 2.025
$arr[0,17] = 3
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = 0.2849999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.825
$arr[0,24] = -1
$arr[0,25] = 1.025
$arr[0,26] = -1
$ws.Range("B144:AB144").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937267
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45350.41666666666
$arr[0,3] = 'Volos NFC'
$arr[0,4] = 'OFI Crete'
$arr[0,5] = 3
$arr[0,6] = 1
$arr[0,7] = 'H'
$arr[0,8] = 2.7
$arr[0,9] = 3.25
$arr[0,10] = 2.625
$arr[0,11] = 2.7
$arr[0,12] = 3.2
$arr[0,13] = 2.8
$arr[0,14] = 0
$arr[0,15] = 1.825
$arr[0,16] = 2.025
$arr[0,17] = 2.25
$arr[0,18] = 2
$arr[0,19] = 1.85
$arr[0,20] = 1.7
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.825
$arr[0,24] = -1
$arr[0,25] = 1
$arr[0,26] = -1
$ws.Range("B168:AB168").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6935703
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45350.41666666666
$arr[0,3] = 'Asteras Tripolis'
$arr[0,4] = 'Kifisias FC'
$arr[0,5] = 3
$arr[0,6] = 3
$arr[0,7] = 'D'
$arr[0,8] = 1.833
$arr[0,9] = 3.4
$arr[0,10] = 4.5
$arr[0,11] = 1.8
$arr[0,12] = 3.5
$arr[0,13] = 4.75
$arr[0,14] = -0.75
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.5
$arr[0,18] = 1.875
$arr[0,19] = 1.975
$arr[0,20] = -1
$arr[0,21] = 2.5
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.825
$arr[0,25] = 0.875
$arr[0,26] = -1
$ws.Range("B169:AB169").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937268
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45350.5
$arr[0,3] = 'Panetolikos'
$arr[0,4] = 'Olympiakos'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 'A'
$arr[0,8] = 8
$arr[0,9] = 5
$arr[0,10] = 1.363
$arr[0,11] = 8.5
$arr[0,12] = 5
$arr[0,13] = 1.363
$arr[0,14] = 1.25
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.75
$arr[0,18] = 1.85
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.363
$arr[0,23] = 0.5125
$arr[0,24] = -0.5
$arr[0,25] = 0.425
$arr[0,26] = -0.5
$ws.Range("B170:AB170").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937266
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45350.5
$arr[0,3] = 'Atromitos Athinon'
$arr[0,4] = 'Lamia'
$arr[0,5] = 3
$arr[0,6] = 1
$arr[0,7] = 'H'
$arr[0,8] = 2.3
$arr[0,9] = 3.2
$arr[0,10] = 3.1
$arr[0,11] = 2.2
$arr[0,12] = 3.3
$arr[0,13] = 3.3
$arr[0,14] = -0.25
$arr[0,15] = 1.925
$arr[0,16] = 1.925
$arr[0,17] = 2.5
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = 1.2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.925
$arr[0,24] = -1
$arr[0,25] = 1.025
$arr[0,26] = -1
$ws.Range("B171:AB171").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6936863
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'OFI Crete'
$arr[0,4] = 'Panathinaikos'
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = 'D'
$arr[0,8] = 8
$arr[0,9] = 4.75
$arr[0,10] = 1.4
$arr[0,11] = 5.5
$arr[0,12] = 4.75
$arr[0,13] = 1.55
$arr[0,14] = 1
$arr[0,15] = 1.95
$arr[0,16] = 1.9
$arr[0,17] = 2.5
$arr[0,18] = 1.85
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = 3.75
$arr[0,22] = -1
$arr[0,23] = 0.95
$arr[0,24] = -1
$arr[0,25] = 0.8500000000000001
$arr[0,26] = -1
$ws.Range("B175:AB175").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937269
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Aris Salonika'
$arr[0,4] = 'AEK Athens'
$arr[0,5] = 3
$arr[0,6] = 3
$arr[0,7] = 'D'
$arr[0,8] = 4.75
$arr[0,9] = 3.75
$arr[0,10] = 1.75
$arr[0,11] = 6.5
$arr[0,12] = 4.2
$arr[0,13] = 1.5
$arr[0,14] = 1
$arr[0,15] = 2.05
$arr[0,16] = 1.8
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.875
$arr[0,20] = -1
$arr[0,21] = 3.2
$arr[0,22] = -1
$arr[0,23] = 1.05
$arr[0,24] = -1
$arr[0,25] = 0.9750000000000001
$arr[0,26] = -1
$ws.Range("B176:AB176").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937270
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Olympiakos'
$arr[0,4] = 'Volos NFC'
$arr[0,5] = 3
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.125
$arr[0,9] = 9
$arr[0,10] = 19
$arr[0,11] = 1.111
$arr[0,12] = 9
$arr[0,13] = 21
$arr[0,14] = -2.25
$arr[0,15] = 1.875
$arr[0,16] = 1.975
$arr[0,17] = 3.25
$arr[0,18] = 2
$arr[0,19] = 1.85
$arr[0,20] = 0.111
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.875
$arr[0,24] = -1
$arr[0,25] = -0.5
$arr[0,26] = 0.425
$ws.Range("B177:AB177").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937271
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Giannina'
$arr[0,4] = 'Atromitos Athinon'
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 'D'
$arr[0,8] = 2.45
$arr[0,9] = 3.1
$arr[0,10] = 3.1
$arr[0,11] = 2
$arr[0,12] = 3.3
$arr[0,13] = 4
$arr[0,14] = -0.5
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.25
$arr[0,18] = 1.85
$arr[0,19] = 2
$arr[0,20] = -1
$arr[0,21] = 2.3
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.825
$arr[0,25] = -0.5
$arr[0,26] = 0.5
$ws.Range("B178:AB178").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6935701
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Kifisias FC'
$arr[0,4] = 'Panetolikos'
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = 'D'
$arr[0,8] = 2.45
$arr[0,9] = 3.25
$arr[0,10] = 3
$arr[0,11] = 2.05
$arr[0,12] = 3.3
$arr[0,13] = 3.8
$arr[0,14] = -0.5
$arr[0,15] = 2.05
$arr[0,16] = 1.8
$arr[0,17] = 2.25
$arr[0,18] = 1.8
$arr[0,19] = 2.05
$arr[0,20] = -1
$arr[0,21] = 2.3
$arr[0,22] = -1
$arr[0,23] = -1
$arr[0,24] = 0.8
$arr[0,25] = 0.8
$arr[0,26] = -1
$ws.Range("B179:AB179").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6937272
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Lamia'
$arr[0,4] = 'PAOK Salonika'
$arr[0,5] = 0
$arr[0,6] = 2
$arr[0,7] = 'A'
$arr[0,8] = 7.5
$arr[0,9] = 4.5
$arr[0,10] = 1.444
$arr[0,11] = 9.5
$arr[0,12] = 5
$arr[0,13] = 1.333
$arr[0,14] = 1.5
$arr[0,15] = 1.925
$arr[0,16] = 1.925
$arr[0,17] = 3
$arr[0,18] = 1.95
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.333
$arr[0,23] = -1
$arr[0,24] = 0.925
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$ws.Range("B180:AB180").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 6935700
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45354.60416666666
$arr[0,3] = 'Panserraikos'
$arr[0,4] = 'Asteras Tripolis'
$arr[0,5] = 2
$arr[0,6] = 1
$arr[0,7] = 'H'
$arr[0,8] = 2.6
$arr[0,9] = 3.2
$arr[0,10] = 2.875
$arr[0,11] = 2.25
$arr[0,12] = 3.3
$arr[0,13] = 3.3
$arr[0,14] = -0.25
$arr[0,15] = 1.925
$arr[0,16] = 1.925
$arr[0,17] = 2.25
$arr[0,18] = 2
$arr[0,19] = 1.85
$arr[0,20] = 1.25
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.925
$arr[0,24] = -1
$arr[0,25] = 1
$arr[0,26] = -1
$ws.Range("B181:AB181").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920471
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45382.60416666666
$arr[0,3] = 'Aris Salonika'
$arr[0,4] = 'Lamia'
$arr[0,5] = 3
$arr[0,6] = 1
$arr[0,7] = 'H'
$arr[0,8] = 1.571
$arr[0,9] = 4
$arr[0,10] = 6
$arr[0,11] = 1.444
$arr[0,12] = 4.5
$arr[0,13] = 8.5
$arr[0,14] = -1.25
$arr[0,15] = 1.925
$arr[0,16] = 1.925
$arr[0,17] = 2.75
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = 0.444
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.925
$arr[0,24] = -1
$arr[0,25] = 1.025
$arr[0,26] = -1
$ws.Range("B194:AB194").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920470
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45382.60416666666
$arr[0,3] = 'AEK Athens'
$arr[0,4] = 'Olympiakos'
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.909
$arr[0,9] = 3.4
$arr[0,10] = 4.2
$arr[0,11] = 2.2
$arr[0,12] = 3.2
$arr[0,13] = 3.5
$arr[0,14] = -0.25
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2.5
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = 1.2
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.8500000000000001
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B195:AB195").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920453
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45388.5625
$arr[0,3] = 'Panetolikos'
$arr[0,4] = 'Volos NFC'
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = 'A'
$arr[0,8] = 2.3
$arr[0,9] = 3
$arr[0,10] = 3.4
$arr[0,11] = 2.1
$arr[0,12] = 3.1
$arr[0,13] = 3.8
$arr[0,14] = -0.25
$arr[0,15] = 1.8
$arr[0,16] = 2.05
$arr[0,17] = 2.25
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.8
$arr[0,23] = -1
$arr[0,24] = 1.05
$arr[0,25] = -1
$arr[0,26] = 0.825
$ws.Range("B200:AB200").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920450
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45388.5625
$arr[0,3] = 'Asteras Tripolis'
$arr[0,4] = 'Kifisias FC'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 'A'
$arr[0,8] = 2.05
$arr[0,9] = 3.3
$arr[0,10] = 3.6
$arr[0,11] = 2.1
$arr[0,12] = 3.5
$arr[0,13] = 3.4
$arr[0,14] = -0.25
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2.75
$arr[0,18] = 2.025
$arr[0,19] = 1.825
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.4
$arr[0,23] = -1
$arr[0,24] = 1
$arr[0,25] = 0.5125
$arr[0,26] = -0.5
$ws.Range("B201:AB201").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920465
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45409.58333333334
$arr[0,3] = 'Panetolikos'
$arr[0,4] = 'Atromitos Athinon'
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.7
$arr[0,9] = 3.6
$arr[0,10] = 5.25
$arr[0,11] = 1.666
$arr[0,12] = 3.75
$arr[0,13] = 5.5
$arr[0,14] = -0.75
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2.5
$arr[0,18] = 1.95
$arr[0,19] = 1.9
$arr[0,20] = 0.6659999999999999
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0.425
$arr[0,24] = -0.5
$arr[0,25] = -1
$arr[0,26] = 0.8999999999999999
$ws.Range("B224:AB224").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920463
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45409.58333333334
$arr[0,3] = 'Volos NFC'
$arr[0,4] = 'Panserraikos'
$arr[0,5] = 1
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.615
$arr[0,9] = 3.9
$arr[0,10] = 5.75
$arr[0,11] = 1.571
$arr[0,12] = 3.6
$arr[0,13] = 7
$arr[0,14] = -1
$arr[0,15] = 1.95
$arr[0,16] = 1.9
$arr[0,17] = 2.5
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = 0.571
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = -1
$arr[0,26] = 0.925
$ws.Range("B225:AB225").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920466
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45423.58333333334
$arr[0,3] = 'Atromitos Athinon'
$arr[0,4] = 'Asteras Tripolis'
$arr[0,5] = 0
$arr[0,6] = 1
$arr[0,7] = 'A'
$arr[0,8] = 2.3
$arr[0,9] = 3.3
$arr[0,10] = 3.2
$arr[0,11] = 2.15
$arr[0,12] = 3.4
$arr[0,13] = 3.3
$arr[0,14] = -0.25
$arr[0,15] = 1.85
$arr[0,16] = 2
$arr[0,17] = 2.5
$arr[0,18] = 1.975
$arr[0,19] = 1.875
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.3
$arr[0,23] = -1
$arr[0,24] = 1
$arr[0,25] = -1
$arr[0,26] = 0.875
$ws.Range("B229:AB229").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920467
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45423.58333333334
$arr[0,3] = 'OFI Crete'
$arr[0,4] = 'Panetolikos'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 'A'
$arr[0,8] = 2
$arr[0,9] = 3.5
$arr[0,10] = 3.75
$arr[0,11] = 2.05
$arr[0,12] = 3.4
$arr[0,13] = 3.5
$arr[0,14] = -0.25
$arr[0,15] = 1.8
$arr[0,16] = 2.05
$arr[0,17] = 2.5
$arr[0,18] = 1.825
$arr[0,19] = 2.025
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 2.5
$arr[0,23] = -1
$arr[0,24] = 1.05
$arr[0,25] = 0.825
$arr[0,26] = -1
$ws.Range("B230:AB230").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920468
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45423.58333333334
$arr[0,3] = 'Panserraikos'
$arr[0,4] = 'Kifisias FC'
$arr[0,5] = 2
$arr[0,6] = 0
$arr[0,7] = 'H'
$arr[0,8] = 1.909
$arr[0,9] = 3.6
$arr[0,10] = 4
$arr[0,11] = 2
$arr[0,12] = 3.6
$arr[0,13] = 3.5
$arr[0,14] = -0.5
$arr[0,15] = 2.025
$arr[0,16] = 1.825
$arr[0,17] = 2.75
$arr[0,18] = 1.825
$arr[0,19] = 2.025
$arr[0,20] = 1
$arr[0,21] = -1
$arr[0,22] = -1
$arr[0,23] = 1.025
$arr[0,24] = -1
$arr[0,25] = -1
$arr[0,26] = 1.025
$ws.Range("B231:AB231").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 7920469
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45423.58333333334
$arr[0,3] = 'Giannina'
$arr[0,4] = 'Volos NFC'
$arr[0,5] = 1
$arr[0,6] = 1
$arr[0,7] = 'D'
$arr[0,8] = 2.375
$arr[0,9] = 3.4
$arr[0,10] = 3
$arr[0,11] = 2.25
$arr[0,12] = 3.3
$arr[0,13] = 3.1
$arr[0,14] = -0.25
$arr[0,15] = 1.95
$arr[0,16] = 1.9
$arr[0,17] = 2.75
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = -1
$arr[0,21] = 2.3
$arr[0,22] = -1
$arr[0,23] = -0.5
$arr[0,24] = 0.45
$arr[0,25] = -1
$arr[0,26] = 0.925
$ws.Range("B232:AB232").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 8140226
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45431.58333333334
$arr[0,3] = 'Aris Salonika'
$arr[0,4] = 'PAOK Salonika'
$arr[0,5] = 1
$arr[0,6] = 2
$arr[0,7] = 'A'
$arr[0,8] = 5.75
$arr[0,9] = 4.5
$arr[0,10] = 1.5
$arr[0,11] = 5.25
$arr[0,12] = 4.1
$arr[0,13] = 1.571
$arr[0,14] = 1
$arr[0,15] = 1.875
$arr[0,16] = 1.975
$arr[0,17] = 2.75
$arr[0,18] = 1.925
$arr[0,19] = 1.925
$arr[0,20] = -1
$arr[0,21] = -1
$arr[0,22] = 0.571
$arr[0,23] = 0
$arr[0,24] = 0
$arr[0,25] = 0.4625
$arr[0,26] = -0.5
$ws.Range("B237:AB237").Value = $arr

$arr = New-Object 'object[,]' 1,27
$arr[0,0] = 8140565
$arr[0,1] = 'Greece Super League 1'
$arr[0,2] = 45431.58333333334
$arr[0,3] = 'Panathinaikos'
$arr[0,4] = 'Olympiakos'
$arr[0,5] = 2
$arr[0,6] = 2
$arr[0,7] = 'D'
$arr[0,8] = 2.4
$arr[0,9] = 3.3
$arr[0,10] = 2.8
$arr[0,11] = 2.875
$arr[0,12] = 3.5
$arr[0,13] = 2.35
$arr[0,14] = 0.25
$arr[0,15] = 1.8
$arr[0,16] = 2.05
$arr[0,17] = 2.75
$arr[0,18] = 1.95
$arr[0,19] = 1.9
$arr[0,20] = -1
$arr[0,21] = 2.5
$arr[0,22] = -1
$arr[0,23] = 0.4
$arr[0,24] = -0.5
$arr[0,25] = 0.95
$arr[0,26] = -1
$ws.Range("B239:AB239").Value = $arr
